# Update the "as_of_utc" timestamp column (AA) for rows 2-26 on both the
# "Главные" and "Линейные" sheets, reflecting the latest stats refresh.
$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-11-03 11:35:25"
$sheetNames = @("Главные", "Линейные")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("AA2:AA26").Value = $newTimestamp
}
